$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "52.026.83"
$ws.Range("E2").Value = "  +0.73%  "

# Row 3
$ws.Range("D3").Value = "2.890.77"
$ws.Range("E3").Value = "  +3.21%  "

# Row 4
$ws.Range("E4").Value = "  +0.05%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "350.53"
$ws.Range("E5").Value = "  -1.43%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "111.25"
$ws.Range("E6").Value = "  +1.53%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.620"
$ws.Range("E9").Value = "  -0.27%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.74"
$ws.Range("E10").Value = "  -1.14%  "

# Row 11
$ws.Range("E11").Value = "  +0.32%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0856"
$ws.Range("E12").Value = "  +1.95%  "

# Row 13
$ws.Range("E13").Value = "  -0.75%  "

# Row 14
$ws.Range("E14").Value = "  -0.50%  "

# Row 15
$ws.Range("D15").Value = "3.344.77"
$ws.Range("E15").Value = "  +3.18%  "

# Row 16
$ws.Range("B16").Value = "Polygon"
$ws.Range("C16").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.996"
$ws.Range("E16").Value = "  +5.63%  "

# Row 17
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "2.887.04"
$ws.Range("E17").Value = "  +3.27%  "

# Row 18
$ws.Range("D18").Value = "52.054.42"
$ws.Range("E18").Value = "  +0.77%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.68"
$ws.Range("E19").Value = "  -1.06%  "

# Row 20
$ws.Range("E20").Value = "  +4.14%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.40"
$ws.Range("E21").Value = "  +7.23%  "

# Row 22
$ws.Range("D22").Value = "0.0₃0978"
$ws.Range("E22").Value = "  +0.61%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.67"
$ws.Range("E23").Value = "  +0.26%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "269.26"
$ws.Range("E24").Value = "  +0.39%  "

# Row 25
$ws.Range("E25").Value = "  -0.09%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.48"
$ws.Range("E26").Value = "  +1.52%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("E27").Value = "  +0.01%  "

# Row 28
$ws.Range("E28").Value = "  -0.59%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.50"
$ws.Range("E29").Value = "  +1.26%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "38.13"
$ws.Range("E30").Value = "  +1.97%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.23"
$ws.Range("E31").Value = "  -0.07%  "

# Row 32
$ws.Range("E32").Value = "  +1.42%  "

# Row 33
$ws.Range("E33").Value = "  +7.79%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0942"
$ws.Range("E34").Value = "  +10.05%  "

# Row 35
$ws.Range("E35").Value = "  +1.28%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0456"
$ws.Range("E36").Value = "  +1.81%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.998"
$ws.Range("E37").Value = "  -0.17%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.28"
$ws.Range("E38").Value = "  +4.41%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.58"

# Row 40
$ws.Range("E40").Value = "  +2.48%  "

# Row 41
$ws.Range("E41").Value = "  +5.78%  "

# Row 42
$ws.Range("E42").Value = "  +1.63%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "22.72"
$ws.Range("E43").Value = "  +2.98%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "122.41"
$ws.Range("E44").Value = "  +1.99%  "

# Row 45
$ws.Range("E45").Value = "  +0.46%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.55"
$ws.Range("E46").Value = "  +3.55%  "

# Row 47
$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.53"
$ws.Range("E47").Value = "  +6.52%  "

# Row 48
$ws.Range("B48").Value = "Maker"
$ws.Range("C48").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D48").Value = "2.194.34"
$ws.Range("E48").Value = "  +2.53%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.270"
$ws.Range("E49").Value = "  +22.55%  "

# Row 50
$ws.Range("E50").Value = "  +1.89%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0322"
$ws.Range("E51").Value = "  +10.16%  "

Write-Host "Done"